$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly Fruta/hortaliza rows (3-6) got their D/M/N/O/P/R/S values
# rotated between rows. New row 3 takes old row 6's values, new row 4
# takes old row 5's, new row 5 takes old row 3's, new row 6 takes old
# row 4's. Apply the resulting (already rotated) target values directly.

$ws.Range("D3").Value = 44320
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 18800
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1044

$ws.Range("D4").Value = 45084
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 972

$ws.Range("D5").Value = 44719
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20400
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1133

$ws.Range("D6").Value = 44362
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 19000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 19500
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 1083
